$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1574.9231
$ws.Range("H12").Value = 2150
$ws.Range("H135").Value = 782.7646999999999
$ws.Range("H139").Value = 69996.336
$ws.Range("H74").Value = 8985.5
$ws.Range("H77").Value = 8985.5
$ws.Range("I12").Value = 400.5
$ws.Range("I135").Value = 788.3570999999999
$ws.Range("I74").Value = 8141.0835
$ws.Range("I77").Value = 8141.0835
$ws.Range("J112").Value = 1748.5
$ws.Range("J12").Value = 2849.8
$ws.Range("J139").Value = 69996.336
$ws.Range("J74").Value = 14052
$ws.Range("J77").Value = 14052
$ws.Range("K12").Value = 400.5
$ws.Range("K135").Value = 7095.2139
$ws.Range("K74").Value = 8141.0835
$ws.Range("K77").Value = 40705.4175
$ws.Range("L112").Value = 5245.5
$ws.Range("L12").Value = 2849.8
$ws.Range("L139").Value = 69996.336
$ws.Range("L74").Value = 14052
$ws.Range("L77").Value = 70260
$ws.Range("M12").Value = -230.5
$ws.Range("M135").Value = -4560.2139
$ws.Range("M74").Value = -7205.0835
$ws.Range("M77").Value = -36025.4175
$ws.Range("N112").Value = -7461.5
$ws.Range("N12").Value = -3189.8
$ws.Range("N139").Value = -80276.336
$ws.Range("N74").Value = -15924
$ws.Range("N77").Value = -79620

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2550
$ws.Range("H45").Value = 58826468
$ws.Range("H5").Value = 337
$ws.Range("H74").Value = 6415829.5
$ws.Range("H77").Value = 6415829.5
$ws.Range("H97").Value = 537.38464
$ws.Range("I35").Value = 2550
$ws.Range("I45").Value = 90910630
$ws.Range("I5").Value = 274.33334
$ws.Range("I74").Value = 9529007
$ws.Range("I77").Value = 9529007
$ws.Range("J35").Value = 0
$ws.Range("J45").Value = 5499
$ws.Range("J5").Value = 399.66666
$ws.Range("J74").Value = 6346.4116
$ws.Range("J77").Value = 6346.4116
$ws.Range("J97").Value = 121.75
$ws.Range("K35").Value = 2550
$ws.Range("K45").Value = 90910630
$ws.Range("K5").Value = 274.33334
$ws.Range("K74").Value = 9529007
$ws.Range("K77").Value = 47645035
$ws.Range("L35").Value = 0
$ws.Range("L45").Value = 5499
$ws.Range("L5").Value = 399.66666
$ws.Range("L74").Value = 6346.4116
$ws.Range("L77").Value = 31732.058
$ws.Range("L97").Value = 121.75
$ws.Range("M35").Value = -2144
$ws.Range("M45").Value = -90910253
$ws.Range("M5").Value = -162.33334
$ws.Range("M74").Value = -9528133
$ws.Range("M77").Value = -47640667
$ws.Range("N35").ClearContents()
$ws.Range("N45").Value = -6253
$ws.Range("N5").Value = -623.66666
$ws.Range("N74").Value = -8094.4116
$ws.Range("N77").Value = -40468.058
$ws.Range("N97").Value = -1113.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 337
$ws.Range("H69").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("H86").Value = 3569.9092
$ws.Range("H89").Value = 3569.9092
$ws.Range("I4").Value = 274.33334
$ws.Range("I86").Value = 1868
$ws.Range("I89").Value = 1868
$ws.Range("J4").Value = 399.66666
$ws.Range("J69").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("J86").Value = 11228.5
$ws.Range("J89").Value = 11228.5
$ws.Range("K4").Value = 274.33334
$ws.Range("K86").Value = 1868
$ws.Range("K89").Value = 9340
$ws.Range("L4").Value = 399.66666
$ws.Range("L69").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("L86").Value = 11228.5
$ws.Range("L89").Value = 56142.5
$ws.Range("M4").Value = -159.33334
$ws.Range("M86").Value = -745
$ws.Range("M89").Value = -3724
$ws.Range("N4").Value = -629.66666
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("N86").Value = -13474.5
$ws.Range("N89").Value = -67374.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 919.6042
$ws.Range("H25").Value = 788.2
$ws.Range("H29").Value = 0
$ws.Range("H8").Value = 701
$ws.Range("H94").Value = 2761.3
$ws.Range("I22").Value = 775.84375
$ws.Range("I25").Value = 451
$ws.Range("I8").Value = 401
$ws.Range("J22").Value = 1207.125
$ws.Range("J25").Value = 1013
$ws.Range("J29").Value = 0
$ws.Range("J8").Value = 1001
$ws.Range("J94").Value = 2948.375
$ws.Range("K22").Value = 775.84375
$ws.Range("K25").Value = 451
$ws.Range("K8").Value = 401
$ws.Range("L22").Value = 1207.125
$ws.Range("L25").Value = 1013
$ws.Range("L29").Value = 0
$ws.Range("L8").Value = 1001
$ws.Range("L94").Value = 2948.375
$ws.Range("M22").Value = -425.84375
$ws.Range("M25").Value = -277
$ws.Range("M8").Value = -261
$ws.Range("N22").Value = -1907.125
$ws.Range("N25").Value = -1361
$ws.Range("N29").Value = 0
$ws.Range("N8").Value = -1281
$ws.Range("N94").Value = -3850.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 150.26666
$ws.Range("H68").Value = 3046
$ws.Range("H71").Value = 3046
$ws.Range("I2").Value = 181
$ws.Range("J68").Value = 3539.7778
$ws.Range("J71").Value = 3539.7778
$ws.Range("K2").Value = 1086
$ws.Range("L68").Value = 10619.3334
$ws.Range("L71").Value = 31858.0002
$ws.Range("M2").Value = -973
$ws.Range("N68").Value = -12241.3334
$ws.Range("N71").Value = -39970.00019999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5227.25
$ws.Range("H113").Value = 5729.8
$ws.Range("H132").Value = 37819.266
$ws.Range("H20").Value = 46000.875
$ws.Range("H80").Value = 6979.9
$ws.Range("H83").Value = 6979.9
$ws.Range("I102").Value = 4218.6665
$ws.Range("I132").Value = 42729.08
$ws.Range("I80").Value = 4828.4287
$ws.Range("I83").Value = 4828.4287
$ws.Range("J113").Value = 5922
$ws.Range("J132").Value = 13270.2
$ws.Range("J20").Value = 49999.5
$ws.Range("J80").Value = 12000
$ws.Range("J83").Value = 12000
$ws.Range("K102").Value = 4218.6665
$ws.Range("K132").Value = 128187.24
$ws.Range("K80").Value = 4828.4287
$ws.Range("K83").Value = 24142.1435
$ws.Range("L113").Value = 5922
$ws.Range("L132").Value = 39810.60000000001
$ws.Range("L20").Value = 49999.5
$ws.Range("L80").Value = 12000
$ws.Range("L83").Value = 60000
$ws.Range("M102").Value = -2596.6665
$ws.Range("M132").Value = -125657.24
$ws.Range("M80").Value = -3830.4287
$ws.Range("M83").Value = -19150.1435
$ws.Range("N113").Value = -10262
$ws.Range("N132").Value = -44870.60000000001
$ws.Range("N20").Value = -50489.5
$ws.Range("N80").Value = -13996
$ws.Range("N83").Value = -69984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5016.32
$ws.Range("H22").Value = 6000
$ws.Range("H27").Value = 6000
$ws.Range("H46").Value = 5583.25
$ws.Range("I136").Value = 2049.9443
$ws.Range("J136").Value = 12644.143
$ws.Range("J22").Value = 6000
$ws.Range("J27").Value = 6000
$ws.Range("J46").Value = 5772.636
$ws.Range("K136").Value = 6149.8329
$ws.Range("L136").Value = 37932.429
$ws.Range("L22").Value = 6000
$ws.Range("L27").Value = 6000
$ws.Range("L46").Value = 5772.636
$ws.Range("M136").Value = -3599.8329
$ws.Range("N136").Value = -43032.429
$ws.Range("N22").Value = -6590
$ws.Range("N27").Value = -6214
$ws.Range("N46").Value = -6148.636

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3407.926
$ws.Range("H126").Value = 3915.2778
$ws.Range("I122").Value = 2680.36
$ws.Range("I126").Value = 3571.6667
$ws.Range("J126").Value = 5633.3335
$ws.Range("K122").Value = 8041.08
$ws.Range("K126").Value = 10715.0001
$ws.Range("L126").Value = 16900.0005
$ws.Range("M122").Value = -5591.08
$ws.Range("M126").Value = -8245.000100000001
$ws.Range("N126").Value = -21840.0005
